# The workbook has a "Tabla1" table (B1:C31) of Experiencia/Salarios values.
# The commit corrects cell C2 ("1.1" years of experience) from "     39,343.00"
# to "     39,343" (drops the redundant ".00").  Re-entering that cell also
# makes Excel re-flow the wrapped-text row heights for the header/first data
# row and for the last block of rows, and moves the active selection to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Core edit: correct the salary text for the "1.1" experience row. The
# original cell used a leading run of non-breaking spaces (as the rest of
# the "Salarios" column does) ahead of the digits, so reproduce that exactly
# and just drop the trailing ".00".
$nbsp = [char]0x00A0
$newSalary = "$nbsp$nbsp$nbsp$nbsp 39,343"
$ws.Range("C2").Value = $newSalary

# Excel reflows the wrap-text row heights for the affected rows after the edit.
$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(2).RowHeight = 18
For ($r = 24; $r -le 31; $r++) {
    $ws.Rows.Item($r).RowHeight = 36
}

# Leave the selection on the edited cell, matching the saved workbook state.
$ws.Range("C2").Select()
